# Auto-generated edit script applying the DM-test results update
# (analisis de las 3 primeras simulaciones)
$wb = $excel.ActiveWorkbook

# --- P_valores sheet: updated p-values for simulations 1-3 (cols E,F and related symmetric cells) ---
$wsP = $wb.Worksheets.Item("P_valores")
$wsP.Range("E2").Value = 0.0002577238745982147
$wsP.Range("F2").Value = 0.00009744071353834372
$wsP.Range("H2").Value = 0.00005522464219898104
$wsP.Range("J2").Value = 0.0001037623964739165
$wsP.Range("E3").Value = 0.001287150121203284
$wsP.Range("F3").Value = 0.0004852385189262343
$wsP.Range("H3").Value = 0.0001243367486258862
$wsP.Range("J3").Value = 0.0001556057810834233
$wsP.Range("E4").Value = 0.714167567523782
$wsP.Range("F4").Value = 0.5022270148748564
$wsP.Range("H4").Value = 0.002682463344695574
$wsP.Range("J4").Value = 0.0005682332554219904
$wsP.Range("B5").Value = 0.0002577238745982147
$wsP.Range("C5").Value = 0.001287150121203284
$wsP.Range("D5").Value = 0.714167567523782
$wsP.Range("F5").Value = 0.3832519255059395
$wsP.Range("G5").Value = 0.1938612494473171
$wsP.Range("H5").Value = 0.00007605091765006122
$wsP.Range("I5").Value = 0.0001435341302695736
$wsP.Range("J5").Value = 0.00008977679979005337
$wsP.Range("B6").Value = 0.00009744071353834372
$wsP.Range("C6").Value = 0.0004852385189262343
$wsP.Range("D6").Value = 0.5022270148748564
$wsP.Range("E6").Value = 0.3832519255059395
$wsP.Range("G6").Value = 0.2723779474529504
$wsP.Range("H6").Value = 0.0003438231333710018
$wsP.Range("I6").Value = 0.0002465096553081469
$wsP.Range("J6").Value = 0.0001425490315443145
$wsP.Range("E7").Value = 0.1938612494473171
$wsP.Range("F7").Value = 0.2723779474529504
$wsP.Range("H7").Value = 0.04963890425688189
$wsP.Range("J7").Value = 0.0000004820615844280951
$wsP.Range("B8").Value = 0.00005522464219898104
$wsP.Range("C8").Value = 0.0001243367486258862
$wsP.Range("D8").Value = 0.002682463344695574
$wsP.Range("E8").Value = 0.00007605091765006122
$wsP.Range("F8").Value = 0.0003438231333710018
$wsP.Range("G8").Value = 0.04963890425688189
$wsP.Range("I8").Value = 0.002386762913121565
$wsP.Range("J8").Value = 0.0006917871841818357
$wsP.Range("E9").Value = 0.0001435341302695736
$wsP.Range("F9").Value = 0.0002465096553081469
$wsP.Range("H9").Value = 0.002386762913121565
$wsP.Range("J9").Value = 0.0001143959611751288
$wsP.Range("B10").Value = 0.0001037623964739165
$wsP.Range("C10").Value = 0.0001556057810834233
$wsP.Range("D10").Value = 0.0005682332554219904
$wsP.Range("E10").Value = 0.00008977679979005337
$wsP.Range("F10").Value = 0.0001425490315443145
$wsP.Range("G10").Value = 0.0000004820615844280951
$wsP.Range("H10").Value = 0.0006917871841818357
$wsP.Range("I10").Value = 0.0001143959611751288

# --- Estadisticos_HLN_DM sheet: updated HLN test statistics ---
$wsE = $wb.Worksheets.Item("Estadisticos_HLN_DM")
$wsE.Range("E2").Value = -4.348212633911547
$wsE.Range("F2").Value = -4.746772442247551
$wsE.Range("H2").Value = -4.980555018623723
$wsE.Range("J2").Value = -4.72095797647928
$wsE.Range("E3").Value = -3.687915806702163
$wsE.Range("F3").Value = -4.089180297696976
$wsE.Range("H3").Value = -4.646729978517622
$wsE.Range("J3").Value = -4.554776532736073
$wsE.Range("E4").Value = -0.3710252963196961
$wsE.Range("F4").Value = -0.6822090562709844
$wsE.Range("H4").Value = -3.382258339709366
$wsE.Range("J4").Value = -4.024456860711632
$wsE.Range("B5").Value = 4.348212633911547
$wsE.Range("C5").Value = 3.687915806702163
$wsE.Range("D5").Value = 0.3710252963196961
$wsE.Range("F5").Value = -0.8897019346994872
$wsE.Range("G5").Value = -1.340191172142515
$wsE.Range("H5").Value = -4.848670148498835
$wsE.Range("I5").Value = -4.587865822655715
$wsE.Range("J5").Value = -4.780430667496843
$wsE.Range("B6").Value = 4.746772442247551
$wsE.Range("C6").Value = 4.089180297696976
$wsE.Range("D6").Value = 0.6822090562709844
$wsE.Range("E6").Value = 0.8897019346994872
$wsE.Range("G6").Value = -1.125809895187498
$wsE.Range("H6").Value = -4.230255406914275
$wsE.Range("I6").Value = -4.36641796718269
$wsE.Range("J6").Value = -4.590688313660952
$wsE.Range("E7").Value = 1.340191172142515
$wsE.Range("F7").Value = 1.125809895187498
$wsE.Range("H7").Value = -2.077452679044677
$wsE.Range("J7").Value = -7.019016154225326
$wsE.Range("B8").Value = 4.980555018623723
$wsE.Range("C8").Value = 4.646729978517622
$wsE.Range("D8").Value = 3.382258339709366
$wsE.Range("E8").Value = 4.848670148498835
$wsE.Range("F8").Value = 4.230255406914275
$wsE.Range("G8").Value = 2.077452679044677
$wsE.Range("I8").Value = -3.431181344536168
$wsE.Range("J8").Value = -3.94371576889124
$wsE.Range("E9").Value = 4.587865822655715
$wsE.Range("F9").Value = 4.36641796718269
$wsE.Range("H9").Value = 3.431181344536168
$wsE.Range("J9").Value = -4.68091333328874
$wsE.Range("B10").Value = 4.72095797647928
$wsE.Range("C10").Value = 4.554776532736073
$wsE.Range("D10").Value = 4.024456860711632
$wsE.Range("E10").Value = 4.780430667496843
$wsE.Range("F10").Value = 4.590688313660952
$wsE.Range("G10").Value = 7.019016154225326
$wsE.Range("H10").Value = 3.94371576889124
$wsE.Range("I10").Value = 4.68091333328874

# --- Resumen_Modelos sheet: updated formatted-text summary values (kept as Text, 4 decimals) ---
$wsR = $wb.Worksheets.Item("Resumen_Modelos")
$wsR.Range("F5:H5").NumberFormat = "@"
$wsR.Range("F6:H6").NumberFormat = "@"
$wsR.Range("F8:H8").NumberFormat = "@"
$wsR.Range("F10:H10").NumberFormat = "@"
$wsR.Range("F5").Value = "3.2318"
$wsR.Range("G5").Value = "5.6931"
$wsR.Range("H5").Value = "1.7616"
$wsR.Range("F6").Value = "3.3418"
$wsR.Range("G6").Value = "5.3088"
$wsR.Range("H6").Value = "1.5886"
$wsR.Range("F8").Value = "5.8505"
$wsR.Range("G8").Value = "8.3502"
$wsR.Range("H8").Value = "1.4273"
$wsR.Range("F10").Value = "11.2519"
$wsR.Range("G10").Value = "19.4286"
$wsR.Range("H10").Value = "1.7267"
